$d = $word.ActiveDocument

function Find-ParaByText($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Range.Text.TrimEnd() -eq $text) {
            return $para
        }
    }
    return $null
}

function Append-Para($afterPara, $text) {
    $afterPara.Range.InsertParagraphAfter()
    $newPara = $afterPara.Next()
    if ($text -ne "") {
        $newPara.Range.Text = $text
    }
    return $newPara
}

# ---------------------------------------------------------------------------
# Hunk 1: the blank separator paragraph right after "Endpoints I/O Domain."
# gets the label "Predictions Domain:" and a new blank paragraph is added
# right after it (so the following "Dimensional Domain, ..." paragraph keeps
# its usual blank-line spacer before it).
# ---------------------------------------------------------------------------
$anchor1 = Find-ParaByText $d "Endpoints I/O Domain."
$p1 = $anchor1.Next()
$p1.Range.Text = "Predictions Domain:"
$p1.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Hunk 2: after the paragraph "OGM / Client Drivers Services." append a whole
# new block of paragraphs, following the document's usual pattern of
# alternating blank separator paragraphs and label/content paragraphs.
# ---------------------------------------------------------------------------
$last = Find-ParaByText $d "OGM / Client Drivers Services."

$last = Append-Para $last ""
$last = Append-Para $last "Predictions Domain:"
$last = Append-Para $last ""
$last = Append-Para $last "Semiotic Layer: (PredictionType, PredictionSubject, PredictionItem, PredictionValue);"
$last = Append-Para $last ""
$last = Append-Para $last "Type Kind: Domain Service Handler. Domain signatures (domain / range: Subject Kind / Object Kind)."
$last = Append-Para $last ""
$last = Append-Para $last "Subject Kind: domain resource types (image)."
$last = Append-Para $last ""
$last = Append-Para $last "Item Kind: prediction resource types (image/face)."
$last = Append-Para $last ""
$last = Append-Para $last "Value Kind: range resource types (face)."
$last = Append-Para $last ""
$last = Append-Para $last "Reify Prediction as Relationship (Values as Relation Resources)."
$last = Append-Para $last ""
$last = Append-Para $last "Domains Dataflow:"
$last = Append-Para $last ""
$last = Append-Para $last "Layers Dataflow: Augmentation. Rules / Productions matching (Reference Model / Kinds Aggregation)."
$last = Append-Para $last ""
$last = Append-Para $last "Semiotic Dataflow: Object Kind matches Subject Kind of Context Kind signatures. (Sucessive Layers Dataflow)."

Write-Output "done"
